# Applies the "Completed all basic functionality" update to the EffectMap
# worksheet: refreshes several precondition/postcondition strings, adds
# EffectID values that were missing, tweaks a couple of probabilities, and
# appends two brand-new effect rows (npc-movement, npc-pick-item).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 2 (room-enter / first visit) - only the postcondition grew a line
# ---------------------------------------------------------------------
$ws.Range("D2").Value = "displayCD(`"room-enter`", {`"roomname`":playerCharacter.getRoom().getRoomName()})`nplayerCharacter.getRoom().getDescription()`nplayerCharacter.getRoom().getItemDescription()`nplayerCharacter.getRoom().getNodeDescription()`nplayerCharacter.getRoom().addVisit()`nnpcintheroom()`nvariableMap[`"JUST_ENTERED`"] =0"

# ---------------------------------------------------------------------
# Row 3 (room-enter / repeat visit) - postcondition rewritten
# ---------------------------------------------------------------------
$ws.Range("D3").Value = "displayCD(`"room-enter`",  {`"roomname`":playerCharacter.getRoom().getRoomName()})`nplayerCharacter.getRoom().getItemDescription()`nplayerCharacter.getRoom().getNodeDescription()`nplayerCharacter.getRoom().addVisit()`nnpcintheroom()`nvariableMap[`"JUST_ENTERED`"] =0"

# ---------------------------------------------------------------------
# Row 4 (game-start) - postcondition rewritten
# ---------------------------------------------------------------------
$ws.Range("D4").Value = "displayCD(`"room-enter`", {`"roomname`":playerCharacter.getRoom().getRoomName()})`nplayerCharacter.getRoom().getDescription()`nplayerCharacter.getRoom().getItemDescription()`nplayerCharacter.getRoom().getNodeDescription()`nnpcintheroom()`nplayerCharacter.getRoom().addVisit()"

# ---------------------------------------------------------------------
# Row 5 (user-death-by-knife) - gained an EffectID, precondition/post
# rewritten to use the npc name + lowercase object + GAME_ACTIVE flag
# ---------------------------------------------------------------------
$ws.Range("A5").Value = 900004
$ws.Range("C5").Value = "npcCharacter.hasObject(`"knife`")`nnpcCharacter.getRoom().getRoomName()==playerCharacter.getRoom().getRoomName()`nnpcCharacter.getAttribute(`"resentment`")>0"
$ws.Range("D5").Value = "displayCD(`"user-death-by-knife`", {`"npcname`":npcCharacter.getName()})`ndisplayCD(`"game-end`", {})`nvariableMap[`"GAME_ACTIVE`"]=False"

# ---------------------------------------------------------------------
# Row 6 (user-death-by-gun)
# ---------------------------------------------------------------------
$ws.Range("A6").Value = 900005
$ws.Range("C6").Value = "npcCharacter.hasObject(`"gun`")`nnpcCharacter.getRoom().getRoomName()==playerCharacter.getRoom().getRoomName()`nnpcCharacter.getAttribute(`"resentment`")>0"
$ws.Range("D6").Value = "displayCD(`"user-death-by-gun`", {`"npcname`":npcCharacter.getName()})`ndisplayCD(`"game-end`", {})`nvariableMap[`"GAME_ACTIVE`"]=False"

# ---------------------------------------------------------------------
# Row 7 (npc-death-by-knife, formerly mislabeled user-suicide text)
# ---------------------------------------------------------------------
$ws.Range("A7").Value = 900006
$ws.Range("C7").Value = "npcCharacter.hasObject(`"knife`")`nnpcCharacter.getRoom().getRoomName()==playerCharacter.getRoom().getRoomName()`nnpcCharacter.getAttribute(`"resentment`")<100"
$ws.Range("D7").Value = "displayCD(`"npc-suicide-by-knife`", {`"npcname`":npcCharacter.getName()})`ndisplayCD(`"game-end`", {})`nvariableMap[`"GAME_ACTIVE`"]=False"

# ---------------------------------------------------------------------
# Row 8 (npc-death-by-gun)
# ---------------------------------------------------------------------
$ws.Range("A8").Value = 900007
$ws.Range("C8").Value = "npcCharacter.hasObject(`"gun`")`nnpcCharacter.getRoom().getRoomName()==playerCharacter.getRoom().getRoomName()`nnpcCharacter.getAttribute(`"resentment`")<100"
$ws.Range("D8").Value = "displayCD(`"npc-suicide-by-gun`", {`"npcname`":npcCharacter.getName()})`ndisplayCD(`"game-end`", {})`nvariableMap[`"GAME_ACTIVE`"]=False"

# ---------------------------------------------------------------------
# Row 9 (npc-puts-batteries-pod) - object name lowercased, probability 80->70
# ---------------------------------------------------------------------
$ws.Range("A9").Value = 900008
$ws.Range("C9").Value = "npcCharacter.hasObject(`"batteries`")`nnpcCharacter.getRoom().getRoomName()==`"Pod room`""
$ws.Range("D9").Value = "variableMap[`"ESCAPEPOD_ACTIVE`"] =True`nnpcCharacter.removeObject(`"batteries`")"
$ws.Range("H9").Value = 70

# ---------------------------------------------------------------------
# Row 10 (npc-escapes-by-pod) - postcondition rewritten, probability 80->90
# ---------------------------------------------------------------------
$ws.Range("A10").Value = 900009
$ws.Range("D10").Value = "displayCD(`"npc-escapes-by-pod`", {`"npcname`":npcCharacter.getName()})`ndisplayCD(`"game-end`", {})`nvariableMap[`"GAME_ACTIVE`"]=False"
$ws.Range("H10").Value = 90

# ---------------------------------------------------------------------
# Row 11 (brand new) - npc-movement
# ---------------------------------------------------------------------
$ws.Range("A11").Value = 900010
$ws.Range("B11").Value = "npc-movement"
$ws.Range("C11").WrapText = $true
$ws.Range("D11").Value = "npcCharacter.move(`"nwse`"[randrange(4)])"
$ws.Range("E11").Value = $false
$ws.Range("F11").Value = -1
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 70

# ---------------------------------------------------------------------
# Row 12 (brand new) - npc-pick-item
# ---------------------------------------------------------------------
$ws.Range("A12").Value = 900011
$ws.Range("B12").Value = "npc-pick-item"
$ws.Range("C12").Value = "len(npcCharacter.getRoom().inv)>0`nlen(npcCharacter.inv)==0"
$ws.Range("D12").Value = "npcCharacter.pick(npcCharacter.getRoom().inv[randrange(len(npcCharacter.getRoom().inv))])`nprint(`"npc has picked`",npcCharacter.inv)"
$ws.Range("E12").Value = $false
$ws.Range("F12").Value = 3
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 90

# ---------------------------------------------------------------------
# Wrap text on the new C/D cells (pre-existing rows already carry the
# wrap-text style; make sure the new ones match too).
# ---------------------------------------------------------------------
$ws.Range("C12").WrapText = $true
$ws.Range("D11").WrapText = $true
$ws.Range("D12").WrapText = $true

# ---------------------------------------------------------------------
# Row heights - Excel recalculated these when the sheet was resaved with
# wrapped multi-line text; set them explicitly to match.
# ---------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 100.8
$ws.Rows.Item(3).RowHeight = 86.4
$ws.Rows.Item(4).RowHeight = 86.4
$ws.Rows.Item(5).RowHeight = 43.2
$ws.Rows.Item(6).RowHeight = 43.2
$ws.Rows.Item(7).RowHeight = 42.45
$ws.Rows.Item(8).RowHeight = 43.2
$ws.Rows.Item(9).RowHeight = 28.8
$ws.Rows.Item(10).RowHeight = 43.2
$ws.Rows.Item(12).RowHeight = 28.8

# ---------------------------------------------------------------------
# Column widths - widened slightly / split C & D into distinct widths.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.11
$ws.Columns.Item(2).ColumnWidth = 20.44
$ws.Columns.Item(3).ColumnWidth = 75.89
$ws.Columns.Item(4).ColumnWidth = 83.44
$ws.Columns.Item(5).ColumnWidth = 19.44
$ws.Columns.Item(6).ColumnWidth = 14.78

# ---------------------------------------------------------------------
# Sheet view tweaks
# ---------------------------------------------------------------------
$ws.Range("H17").Select()
